$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 491.33334
$ws.Range("I33").Value = 172.18182
$ws.Range("J33").Value = 4002
$ws.Range("K33").Value = 172.18182
$ws.Range("L33").Value = 4002
$ws.Range("M33").Value = 56.81818000000001
$ws.Range("N33").Value = -4460

$ws.Range("H58").Value = 1634.3636
$ws.Range("I58").Value = 884.875
$ws.Range("J58").Value = 3633
$ws.Range("K58").Value = 2654.625
$ws.Range("L58").Value = 10899
$ws.Range("M58").Value = -2504.625

$ws.Range("H62").Value = 35728468
$ws.Range("I62").Value = 66675004
$ws.Range("J62").Value = 20922.154
$ws.Range("K62").Value = 66675004
$ws.Range("L62").Value = 20922.154
$ws.Range("M62").Value = -66674380

$ws.Range("H65").Value = 35728468
$ws.Range("I65").Value = 66675004
$ws.Range("J65").Value = 20922.154
$ws.Range("K65").Value = 333375020
$ws.Range("L65").Value = 104610.77
$ws.Range("M65").Value = -333371900

$ws.Range("H98").Value = 33972.777
$ws.Range("I98").Value = 33030
$ws.Range("J98").Value = 50000
$ws.Range("K98").Value = 33030
$ws.Range("L98").Value = 50000
$ws.Range("M98").Value = -31532

$ws.Range("H113").Value = 22449
$ws.Range("I113").Value = 22449
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 22449
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -19195
$ws.Range("N113").ClearContents()

$ws.Range("H116").Value = 621956.6
$ws.Range("I116").Value = 1392747.2
$ws.Range("J116").Value = 5324.1
$ws.Range("K116").Value = 1392747.2
$ws.Range("L116").Value = 5324.1
$ws.Range("M116").Value = -1389305.2
$ws.Range("N116").Value = -12208.1

$ws.Range("H122").Value = 33972.777
$ws.Range("I122").Value = 33030
$ws.Range("J122").Value = 50000
$ws.Range("K122").Value = 99090
$ws.Range("L122").Value = 150000
$ws.Range("M122").Value = -96640

$ws.Range("H127").Value = 1079.7
$ws.Range("I127").Value = 921.8889
$ws.Range("J127").Value = 2500
$ws.Range("K127").Value = 2765.6667
$ws.Range("L127").Value = 7500
$ws.Range("M127").Value = 2194.3333

$ws.Range("H131").Value = 6951.9062
$ws.Range("I131").Value = 1871.6666
$ws.Range("J131").Value = 8939.825999999999
$ws.Range("K131").Value = 5614.9998
$ws.Range("L131").Value = 26819.478
$ws.Range("M131").Value = -574.9997999999996
$ws.Range("N131").Value = -36899.478

$ws.Range("H137").Value = 465237.03
$ws.Range("I137").Value = 585786.5600000001
$ws.Range("J137").Value = 31258.8
$ws.Range("K137").Value = 1757359.68
$ws.Range("L137").Value = 93776.39999999999
$ws.Range("M137").Value = -1754809.68

$ws.Range("H138").Value = 4445.1323
$ws.Range("I138").Value = 1108.95
$ws.Range("J138").Value = 5835.2085
$ws.Range("K138").Value = 3326.85
$ws.Range("L138").Value = 17505.6255
$ws.Range("M138").Value = 1813.15
$ws.Range("N138").Value = -27785.6255

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4705.963
$ws.Range("I61").Value = 4797.9614
$ws.Range("J61").Value = 2314
$ws.Range("K61").Value = 4797.9614
$ws.Range("L61").Value = 2314
$ws.Range("M61").Value = -4585.9614

$ws.Range("H97").Value = 49685.145
$ws.Range("I97").Value = 22517
$ws.Range("J97").Value = 117605.5
$ws.Range("K97").Value = 22517
$ws.Range("L97").Value = 117605.5
$ws.Range("M97").Value = -22021
$ws.Range("N97").Value = -118597.5

$ws.Range("H132").Value = 2732.2727
$ws.Range("I132").Value = 2246.72
$ws.Range("J132").Value = 4249.625
$ws.Range("K132").Value = 6740.16
$ws.Range("L132").Value = 12748.875
$ws.Range("M132").Value = -4210.16

$ws.Range("H136").Value = 4705.963
$ws.Range("I136").Value = 4797.9614
$ws.Range("J136").Value = 2314
$ws.Range("K136").Value = 14393.8842
$ws.Range("L136").Value = 6942
$ws.Range("M136").Value = -11843.8842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 65940.25
$ws.Range("I82").Value = 16664
$ws.Range("J82").Value = 95506
$ws.Range("K82").Value = 16664
$ws.Range("L82").Value = 95506
$ws.Range("M82").Value = -16281

$ws.Range("H85").Value = 65940.25
$ws.Range("I85").Value = 16664
$ws.Range("J85").Value = 95506
$ws.Range("K85").Value = 16664
$ws.Range("L85").Value = 95506
$ws.Range("M85").Value = -15338

$ws.Range("H94").Value = 1645.1666
$ws.Range("I94").Value = 1229.5
$ws.Range("J94").Value = 2476.5
$ws.Range("K94").Value = 1229.5
$ws.Range("L94").Value = 2476.5
$ws.Range("M94").Value = -778.5

$ws.Range("H134").Value = 7470.6294
$ws.Range("I134").Value = 9024.263000000001
$ws.Range("J134").Value = 3780.75
$ws.Range("K134").Value = 27072.789
$ws.Range("L134").Value = 11342.25
$ws.Range("M134").Value = -24537.789

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 11750
$ws.Range("I41").Value = 11750
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 11750
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -11322
$ws.Range("N41").ClearContents()

$ws.Range("H51").Value = 34999
$ws.Range("I51").Value = 34999
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 34999
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -34263
$ws.Range("N51").ClearContents()

$ws.Range("H61").Value = 34999
$ws.Range("I61").Value = 34999
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 34999
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -34651
$ws.Range("N61").ClearContents()

$ws.Range("H99").Value = 23229422
$ws.Range("I99").Value = 23229422
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 23229422
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -23227924

$ws.Range("H122").Value = 2306.75
$ws.Range("I122").Value = 2188.6
$ws.Range("J122").Value = 2503.6667
$ws.Range("K122").Value = 6565.799999999999
$ws.Range("L122").Value = 7511.000100000001
$ws.Range("M122").Value = -4115.799999999999
$ws.Range("N122").Value = -12411.0001

$ws.Range("H126").Value = 23229422
$ws.Range("I126").Value = 23229422
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 69688266
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -69685796

$ws.Range("H132").Value = 50117.91
$ws.Range("I132").Value = 16588.555
$ws.Range("J132").Value = 201000
$ws.Range("K132").Value = 49765.665
$ws.Range("L132").Value = 603000
$ws.Range("M132").Value = -47235.665
$ws.Range("N132").Value = -608060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 2066.3333
$ws.Range("I123").Value = 2066.3333
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 6198.999899999999
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -3748.999899999999
$ws.Range("N123").ClearContents()

$ws.Range("H137").Value = 5205.8
$ws.Range("I137").Value = 3015
$ws.Range("J137").Value = 6666.3335
$ws.Range("K137").Value = 9045
$ws.Range("L137").Value = 19999.0005
$ws.Range("M137").Value = -3945
$ws.Range("N137").Value = -30199.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 85.52631
$ws.Range("I2").Value = 62.5
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 62.5
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = 50.5

$ws.Range("H7").Value = 26666666
$ws.Range("I7").Value = 35000000
$ws.Range("J7").Value = 10000000
$ws.Range("K7").Value = 35000000
$ws.Range("L7").Value = 10000000
$ws.Range("M7").Value = -34999888

$ws.Range("H8").Value = 26666666
$ws.Range("I8").Value = 35000000
$ws.Range("J8").Value = 10000000
$ws.Range("K8").Value = 35000000
$ws.Range("L8").Value = 10000000
$ws.Range("M8").Value = -34999861

$ws.Range("H102").Value = 10802
$ws.Range("I102").Value = 10802
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 10802
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -9180

$ws.Range("H122").Value = 14356.619
$ws.Range("I122").Value = 16939.938
$ws.Range("J122").Value = 6090
$ws.Range("K122").Value = 50819.814
$ws.Range("L122").Value = 18270
$ws.Range("M122").Value = -48369.814
$ws.Range("N122").Value = -23170

$ws.Range("H126").Value = 17252.916
$ws.Range("I126").Value = 18064.908
$ws.Range("J126").Value = 16565.846
$ws.Range("K126").Value = 54194.724
$ws.Range("L126").Value = 49697.538
$ws.Range("M126").Value = -51724.724

$ws.Range("H138").Value = 59000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 59000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 59000
$ws.Range("N138").Value = -69280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2321.1853
$ws.Range("I46").Value = 1323.875
$ws.Range("J46").Value = 2741.1052
$ws.Range("K46").Value = 1323.875
$ws.Range("L46").Value = 2741.1052
$ws.Range("M46").Value = -1135.875
$ws.Range("N46").Value = -3117.1052

$ws.Range("H100").Value = 5424.8237
$ws.Range("I100").Value = 6247.6665
$ws.Range("J100").Value = 3450
$ws.Range("K100").Value = 6247.6665
$ws.Range("L100").Value = 3450
$ws.Range("M100").Value = -5706.6665
$ws.Range("N100").Value = -4532

$ws.Range("H122").Value = 4672.579
$ws.Range("I122").Value = 4480.364
$ws.Range("J122").Value = 4936.875
$ws.Range("K122").Value = 13441.092
$ws.Range("L122").Value = 14810.625
$ws.Range("M122").Value = -10991.092
$ws.Range("N122").Value = -19710.625

$ws.Range("H132").Value = 1366848.8
$ws.Range("I132").Value = 3000346.5
$ws.Range("J132").Value = 5600.6665
$ws.Range("K132").Value = 9001039.5
$ws.Range("L132").Value = 16801.9995
$ws.Range("M132").Value = -8998509.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

$ws.Range("H10").Value = 1201.3334
$ws.Range("I10").Value = 105
$ws.Range("J10").Value = 1749.5
$ws.Range("K10").Value = 105
$ws.Range("L10").Value = 1749.5
$ws.Range("M10").Value = 64
$ws.Range("N10").Value = -2087.5

$ws.Range("H45").Value = 18965
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 18965
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 18965
$ws.Range("N45").Value = -19947

$ws.Range("H62").Value = 85320.56
$ws.Range("I62").Value = 156234
$ws.Range("J62").Value = 3210.2632
$ws.Range("K62").Value = 156234
$ws.Range("L62").Value = 3210.2632
$ws.Range("M62").Value = -155610
$ws.Range("N62").Value = -4458.263199999999

$ws.Range("H65").Value = 85320.56
$ws.Range("I65").Value = 156234
$ws.Range("J65").Value = 3210.2632
$ws.Range("K65").Value = 781170
$ws.Range("L65").Value = 16051.316
$ws.Range("M65").Value = -778050
$ws.Range("N65").Value = -22291.316

$ws.Range("H81").Value = 27998.182
$ws.Range("I81").Value = 30560.125
$ws.Range("J81").Value = 21166.334
$ws.Range("K81").Value = 61120.25
$ws.Range("L81").Value = 42332.668
$ws.Range("M81").Value = -60059.25

$ws.Range("H84").Value = 27998.182
$ws.Range("I84").Value = 30560.125
$ws.Range("J84").Value = 21166.334
$ws.Range("K84").Value = 305601.25
$ws.Range("L84").Value = 211663.34
$ws.Range("M84").Value = -300297.25

$ws.Range("H100").Value = 32318.059
$ws.Range("I100").Value = 17329.143
$ws.Range("J100").Value = 102266.336
$ws.Range("K100").Value = 34658.286
$ws.Range("L100").Value = 204532.672
$ws.Range("M100").Value = -34117.286

$ws.Range("H122").Value = 2486.1538
$ws.Range("I122").Value = 1893
$ws.Range("J122").Value = 5748.5
$ws.Range("K122").Value = 5679
$ws.Range("L122").Value = 17245.5
$ws.Range("M122").Value = -3229
$ws.Range("N122").Value = -22145.5

$ws.Range("H132").Value = 18920.5
$ws.Range("I132").Value = 21507.814
$ws.Range("J132").Value = 8940.857
$ws.Range("K132").Value = 64523.442
$ws.Range("L132").Value = 26822.571
$ws.Range("M132").Value = -61993.442
$ws.Range("N132").Value = -31882.571
